$d = $word.ActiveDocument
$d.Content.Find.Execute("31+1=32", $true, $true, $false, $false, $false, $true, 1, $false, "27+45=72", 2)
$d.Content.Find.Execute("66+13=79", $true, $true, $false, $false, $false, $true, 1, $false, "83-24=59", 2)
$d.Content.Find.Execute("61+7=68", $true, $true, $false, $false, $false, $true, 1, $false, "58-26=32", 2)
$d.Content.Find.Execute("5+39=44", $true, $true, $false, $false, $false, $true, 1, $false, "94-88=6", 2)
$d.Content.Find.Execute("16+66=82", $true, $true, $false, $false, $false, $true, 1, $false, "80-24=56", 2)
$d.Content.Find.Execute("24-17=7", $true, $true, $false, $false, $false, $true, 1, $false, "98-88=10", 2)
$d.Content.Find.Execute("96-80=16", $true, $true, $false, $false, $false, $true, 1, $false, "27+69=96", 2)
$d.Content.Find.Execute("85-5=80", $true, $true, $false, $false, $false, $true, 1, $false, "70+3=73", 2)
$d.Content.Find.Execute("68-1=67", $true, $true, $false, $false, $false, $true, 1, $false, "39+48=87", 2)
$d.Content.Find.Execute("22+58=80", $true, $true, $false, $false, $false, $true, 1, $false, "33-26=7", 2)
$d.Content.Find.Execute("84-51=33", $true, $true, $false, $false, $false, $true, 1, $false, "78-12=66", 2)
$d.Content.Find.Execute("11+77=88", $true, $true, $false, $false, $false, $true, 1, $false, "62+1=63", 2)
$d.Content.Find.Execute("43+19=62", $true, $true, $false, $false, $false, $true, 1, $false, "56-18=38", 2)
$d.Content.Find.Execute("13+31=44", $true, $true, $false, $false, $false, $true, 1, $false, "10-9=1", 2)
$d.Content.Find.Execute("86-63=23", $true, $true, $false, $false, $false, $true, 1, $false, "43-29=14", 2)
$d.Content.Find.Execute("56-35=21", $true, $true, $false, $false, $false, $true, 1, $false, "56-20=36", 2)
$d.Content.Find.Execute("90-4=86", $true, $true, $false, $false, $false, $true, 1, $false, "10-7=3", 2)
$d.Content.Find.Execute("65+11=76", $true, $true, $false, $false, $false, $true, 1, $false, "91-90=1", 2)
$d.Content.Find.Execute("95-31=64", $true, $true, $false, $false, $false, $true, 1, $false, "60+8=68", 2)
$d.Content.Find.Execute("15-1=14", $true, $true, $false, $false, $false, $true, 1, $false, "38+53=91", 2)
$d.Content.Find.Execute("20+11=31", $true, $true, $false, $false, $false, $true, 1, $false, "47+23=70", 2)
$d.Content.Find.Execute("21-7=14", $true, $true, $false, $false, $false, $true, 1, $false, "47+3=50", 2)
$d.Content.Find.Execute("99-65=34", $true, $true, $false, $false, $false, $true, 1, $false, "41-22=19", 2)
$d.Content.Find.Execute("57+7=64", $true, $true, $false, $false, $false, $true, 1, $false, "43-28=15", 2)
$d.Content.Find.Execute("84+13=97", $true, $true, $false, $false, $false, $true, 1, $false, "26-0=26", 2)
$d.Content.Find.Execute("44-30=14", $true, $true, $false, $false, $false, $true, 1, $false, "87-20=67", 2)
$d.Content.Find.Execute("97-28=69", $true, $true, $false, $false, $false, $true, 1, $false, "49-22=27", 2)
$d.Content.Find.Execute("62-11=51", $true, $true, $false, $false, $false, $true, 1, $false, "52-0=52", 2)
$d.Content.Find.Execute("46+8=54", $true, $true, $false, $false, $false, $true, 1, $false, "28-0=28", 2)
$d.Content.Find.Execute("98-70=28", $true, $true, $false, $false, $false, $true, 1, $false, "67-29=38", 2)
$d.Content.Find.Execute("36+14=50", $true, $true, $false, $false, $false, $true, 1, $false, "93-85=8", 2)
$d.Content.Find.Execute("72+13=85", $true, $true, $false, $false, $false, $true, 1, $false, "27+26=53", 2)
$d.Content.Find.Execute("44+6=50", $true, $true, $false, $false, $false, $true, 1, $false, "56+18=74", 2)
$d.Content.Find.Execute("74-16=58", $true, $true, $false, $false, $false, $true, 1, $false, "13+81=94", 2)
$d.Content.Find.Execute("84-8=76", $true, $true, $false, $false, $false, $true, 1, $false, "34-21=13", 2)
$d.Content.Find.Execute("21+33=54", $true, $true, $false, $false, $false, $true, 1, $false, "67-41=26", 2)
$d.Content.Find.Execute("70-39=31", $true, $true, $false, $false, $false, $true, 1, $false, "52-12=40", 2)
$d.Content.Find.Execute("45+36=81", $true, $true, $false, $false, $false, $true, 1, $false, "27-20=7", 2)
$d.Content.Find.Execute("67+32=99", $true, $true, $false, $false, $false, $true, 1, $false, "76-65=11", 2)
$d.Content.Find.Execute("77-62=15", $true, $true, $false, $false, $false, $true, 1, $false, "71-29=42", 2)
$d.Content.Find.Execute("98-60=38", $true, $true, $false, $false, $false, $true, 1, $false, "82-59=23", 2)
$d.Content.Find.Execute("88-51=37", $true, $true, $false, $false, $false, $true, 1, $false, "3+71=74", 2)
$d.Content.Find.Execute("20+77=97", $true, $true, $false, $false, $false, $true, 1, $false, "33+51=84", 2)
$d.Content.Find.Execute("57+10=67", $true, $true, $false, $false, $false, $true, 1, $false, "67+11=78", 2)
$d.Content.Find.Execute("75-32=43", $true, $true, $false, $false, $false, $true, 1, $false, "51-0=51", 2)
$d.Content.Find.Execute("12+40=52", $true, $true, $false, $false, $false, $true, 1, $false, "21+22=43", 2)
$d.Content.Find.Execute("6-5=1", $true, $true, $false, $false, $false, $true, 1, $false, "39+49=88", 2)
$d.Content.Find.Execute("33+21=54", $true, $true, $false, $false, $false, $true, 1, $false, "85-84=1", 2)
$d.Content.Find.Execute("40-23=17", $true, $true, $false, $false, $false, $true, 1, $false, "27+42=69", 2)
$d.Content.Find.Execute("95-65=30", $true, $true, $false, $false, $false, $true, 1, $false, "6+53=59", 2)
$d.Content.Find.Execute("50-48=2", $true, $true, $false, $false, $false, $true, 1, $false, "59-36=23", 2)
$d.Content.Find.Execute("2+27=29", $true, $true, $false, $false, $false, $true, 1, $false, "55-16=39", 2)
$d.Content.Find.Execute("9+20=29", $true, $true, $false, $false, $false, $true, 1, $false, "77+8=85", 2)
$d.Content.Find.Execute("76-70=6", $true, $true, $false, $false, $false, $true, 1, $false, "45+44=89", 2)
$d.Content.Find.Execute("14+64=78", $true, $true, $false, $false, $false, $true, 1, $false, "71-41=30", 2)
$d.Content.Find.Execute("46+5=51", $true, $true, $false, $false, $false, $true, 1, $false, "84-64=20", 2)
$d.Content.Find.Execute("99-24=75", $true, $true, $false, $false, $false, $true, 1, $false, "89-52=37", 2)
$d.Content.Find.Execute("83-29=54", $true, $true, $false, $false, $false, $true, 1, $false, "39+11=50", 2)
$d.Content.Find.Execute("52-13=39", $true, $true, $false, $false, $false, $true, 1, $false, "48+6=54", 2)
$d.Content.Find.Execute("32+59=91", $true, $true, $false, $false, $false, $true, 1, $false, "29+15=44", 2)
$d.Content.Find.Execute("11+21=32", $true, $true, $false, $false, $false, $true, 1, $false, "46-6=40", 2)
$d.Content.Find.Execute("10+15=25", $true, $true, $false, $false, $false, $true, 1, $false, "14+5=19", 2)
$d.Content.Find.Execute("47-33=14", $true, $true, $false, $false, $false, $true, 1, $false, "14+21=35", 2)
$d.Content.Find.Execute("79-40=39", $true, $true, $false, $false, $false, $true, 1, $false, "25+1=26", 2)
$d.Content.Find.Execute("53+5=58", $true, $true, $false, $false, $false, $true, 1, $false, "90-53=37", 2)
$d.Content.Find.Execute("40+21=61", $true, $true, $false, $false, $false, $true, 1, $false, "60-41=19", 2)
$d.Content.Find.Execute("8+21=29", $true, $true, $false, $false, $false, $true, 1, $false, "44+53=97", 2)
$d.Content.Find.Execute("84-20=64", $true, $true, $false, $false, $false, $true, 1, $false, "73-67=6", 2)
$d.Content.Find.Execute("93-6=87", $true, $true, $false, $false, $false, $true, 1, $false, "90-43=47", 2)
$d.Content.Find.Execute("95-83=12", $true, $true, $false, $false, $false, $true, 1, $false, "94+1=95", 2)
$d.Content.Find.Execute("12+58=70", $true, $true, $false, $false, $false, $true, 1, $false, "58-10=48", 2)
$d.Content.Find.Execute("55-17=38", $true, $true, $false, $false, $false, $true, 1, $false, "21+59=80", 2)
$d.Content.Find.Execute("15+21=36", $true, $true, $false, $false, $false, $true, 1, $false, "95-15=80", 2)
$d.Content.Find.Execute("72-66=6", $true, $true, $false, $false, $false, $true, 1, $false, "95-36=59", 2)
$d.Content.Find.Execute("52-36=16", $true, $true, $false, $false, $false, $true, 1, $false, "8+69=77", 2)
$d.Content.Find.Execute("53+36=89", $true, $true, $false, $false, $false, $true, 1, $false, "4+90=94", 2)
$d.Content.Find.Execute("57-9=48", $true, $true, $false, $false, $false, $true, 1, $false, "36+21=57", 2)
$d.Content.Find.Execute("4+26=30", $true, $true, $false, $false, $false, $true, 1, $false, "74-37=37", 2)
$d.Content.Find.Execute("56-21=35", $true, $true, $false, $false, $false, $true, 1, $false, "69-66=3", 2)
$d.Content.Find.Execute("70-61=9", $true, $true, $false, $false, $false, $true, 1, $false, "88-2=86", 2)
$d.Content.Find.Execute("12+50=62", $true, $true, $false, $false, $false, $true, 1, $false, "95-16=79", 2)
$d.Content.Find.Execute("21+57=78", $true, $true, $false, $false, $false, $true, 1, $false, "81-3=78", 2)
$d.Content.Find.Execute("78-1=77", $true, $true, $false, $false, $false, $true, 1, $false, "97-70=27", 2)
$d.Content.Find.Execute("39+38=77", $true, $true, $false, $false, $false, $true, 1, $false, "59+11=70", 2)
$d.Content.Find.Execute("7-4=3", $true, $true, $false, $false, $false, $true, 1, $false, "9+89=98", 2)
$d.Content.Find.Execute("32+17=49", $true, $true, $false, $false, $false, $true, 1, $false, "98-53=45", 2)
$d.Content.Find.Execute("99-21=78", $true, $true, $false, $false, $false, $true, 1, $false, "72-37=35", 2)
$d.Content.Find.Execute("94-93=1", $true, $true, $false, $false, $false, $true, 1, $false, "45+41=86", 2)
$d.Content.Find.Execute("76+9=85", $true, $true, $false, $false, $false, $true, 1, $false, "4+80=84", 2)
$d.Content.Find.Execute("45-23=22", $true, $true, $false, $false, $false, $true, 1, $false, "25-14=11", 2)
$d.Content.Find.Execute("62-47=15", $true, $true, $false, $false, $false, $true, 1, $false, "5+63=68", 2)
$d.Content.Find.Execute("69-57=12", $true, $true, $false, $false, $false, $true, 1, $false, "32-6=26", 2)
$d.Content.Find.Execute("38-10=28", $true, $true, $false, $false, $false, $true, 1, $false, "91-23=68", 2)
$d.Content.Find.Execute("30+20=50", $true, $true, $false, $false, $false, $true, 1, $false, "49-27=22", 2)
$d.Content.Find.Execute("34+52=86", $true, $true, $false, $false, $false, $true, 1, $false, "29+49=78", 2)
$d.Content.Find.Execute("18+71=89", $true, $true, $false, $false, $false, $true, 1, $false, "95-58=37", 2)
$d.Content.Find.Execute("48-19=29", $true, $true, $false, $false, $false, $true, 1, $false, "65-16=49", 2)
$d.Content.Find.Execute("74+7=81", $true, $true, $false, $false, $false, $true, 1, $false, "47-41=6", 2)
$d.Content.Find.Execute("99-8=91", $true, $true, $false, $false, $false, $true, 1, $false, "99-23=76", 2)
$d.Content.Find.Execute("40-27=13", $true, $true, $false, $false, $false, $true, 1, $false, "53-19=34", 2)
